# Applies the cryptos price/volume update as captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.694.61"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "'3.793.51"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'432.02"
$ws.Range("E5").Value = "  +4.62%  "
$ws.Range("D6").Value = "'141.06"
$ws.Range("E6").Value = "  +6.60%  "
$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -11.48%  "
$ws.Range("D11").Value = "'0.0000315"
$ws.Range("E11").Value = "  -17.05%  "
$ws.Range("D12").Value = "'42.94"
$ws.Range("E12").Value = "  +4.46%  "
$ws.Range("D13").Value = "'10.47"
$ws.Range("E13").Value = "  +4.08%  "
$ws.Range("D14").Value = "'4.391.41"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "'14.92"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "'3.817.65"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'19.95"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("E19").Value = "  +6.16%  "
$ws.Range("D20").Value = "'66.799.90"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "'410.03"
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").Value = "'14.77"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  +5.73%  "
$ws.Range("D24").Value = "'85.38"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("D25").Value = "'36.82"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +6.61%  "
$ws.Range("D27").Value = "'5.63"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "'9.57"
$ws.Range("E28").Value = "  +32.28%  "
$ws.Range("D29").Value = "'9.79"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").Value = "'722.29"
$ws.Range("E30").Value = "  +5.03%  "
$ws.Range("D31").Value = "'13.79"
$ws.Range("E31").Value = "  +9.93%  "
$ws.Range("E32").Value = "  +10.48%  "
$ws.Range("D33").Value = "'2.69"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'41.76"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'5.71"
$ws.Range("E36").Value = "  +27.33%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "'55.97"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("D39").Value = "'0.0475"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "  +42.88%  "
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -4.99%  "

# Rows 42 and 43 swap positions (PEPE <-> Stellar) with updated price/volume
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.141"
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "'0.0₃0678"
$ws.Range("E43").Value = "  -17.30%  "

$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "'3.25"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'0.321"
$ws.Range("E47").Value = "  +8.54%  "
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("D49").Value = "'2.09"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "'142.44"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("E51").Value = "  -0.09%  "
